$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the year header row (row 4) with 2021 (M) and 2022 (N), ---
# --- matching the existing "year" cell style (s="13", copied from L4). ---
$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4:N4").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").Value = 2021
$ws.Range("N4").Value = 2022

# --- Row 5 (style s="19", copied from L5) ---
$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5:N5").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Value = 5.6
$ws.Range("N5").Value = 6.3

# --- Row 6 (style s="20", copied from L6) ---
$ws.Range("L6").Copy() | Out-Null
$ws.Range("M6:N6").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").Value = 0.8
$ws.Range("N6").Value = 0.8

# --- Row 7 (style s="20", copied from L7) ---
$ws.Range("L7").Copy() | Out-Null
$ws.Range("M7:N7").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Value = 1.9
$ws.Range("N7").Value = 2.4

# --- Row 8 (style s="20", copied from L8) ---
$ws.Range("L8").Copy() | Out-Null
$ws.Range("M8:N8").PasteSpecial(-4122) | Out-Null
$ws.Range("M8").Value = 0.7
$ws.Range("N8").Value = 0.7

# --- Row 9 (style s="20", copied from L9) ---
$ws.Range("L9").Copy() | Out-Null
$ws.Range("M9:N9").PasteSpecial(-4122) | Out-Null
$ws.Range("M9").Value = 0.7
$ws.Range("N9").Value = 0.8

# --- Row 10 (style s="20", copied from L10) ---
$ws.Range("L10").Copy() | Out-Null
$ws.Range("M10:N10").PasteSpecial(-4122) | Out-Null
$ws.Range("M10").Value = 0.9
$ws.Range("N10").Value = 1

# --- Row 11 (style s="20", copied from L11) ---
$ws.Range("L11").Copy() | Out-Null
$ws.Range("M11:N11").PasteSpecial(-4122) | Out-Null
$ws.Range("M11").Value = 0.3
$ws.Range("N11").Value = 0.2

# --- Row 12 (style s="21", copied from L12) ---
$ws.Range("L12").Copy() | Out-Null
$ws.Range("M12:N12").PasteSpecial(-4122) | Out-Null
$ws.Range("M12").Value = 0.2
$ws.Range("N12").Value = 0.4

# --- New row 14, footnote about the 2022 forestry-service data source. ---
# B13 carries the existing footnote style (s="17") - copy its formatting
# onto B14 and fill in the new text there.
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$ws.Rows.Item(14).RowHeight = 34.5

$excel.CutCopyMode = $false
